# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$Cell, [string]$Text)
    $range = $Sheet.Range($Cell)
    # Force text storage so numeric-looking strings (e.g. "212.74")
    # are kept verbatim instead of being coerced into numbers.
    $range.NumberFormat = "@"
    $range.Value = $Text
}

# --- Rows 28/29 swapped places (BinanceUSD now ranks above EthereumClassic) ---
$ws.Range("B28").Value = 'BinanceUSD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws "D28" '1.00'
$ws.Range("E28").Value = '  -0.19%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws "D29" '15.65'
$ws.Range("E29").Value = '  -4.41%  '

# --- Refreshed price / volume figures for the remaining coins ---
$ws.Range("D2").Value = '27.515.65'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '1.646.54'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  -0.09%  '
Set-TextValue $ws "D5" '212.74'
$ws.Range("E5").Value = '  -1.33%  '
Set-TextValue $ws "D6" '0.529'
$ws.Range("E6").Value = '  +3.96%  '
Set-TextValue $ws "D8" '23.54'
$ws.Range("E8").Value = '  -2.26%  '
Set-TextValue $ws "D9" '0.257'
$ws.Range("E10").Value = '  -1.37%  '
Set-TextValue $ws "D11" '0.0893'
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("D12").Value = '1.878.99'
$ws.Range("E12").Value = '  -1.23%  '
$ws.Range("D13").Value = '1.646.07'
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("E14").Value = '  +3.76%  '
$ws.Range("E15").Value = '  -2.11%  '
Set-TextValue $ws "D16" '64.53'
$ws.Range("E16").Value = '  -2.89%  '
$ws.Range("D17").Value = '27.480.89'
$ws.Range("E17").Value = '  -0.46%  '
Set-TextValue $ws "D18" '231.60'
$ws.Range("E18").Value = '  -4.25%  '
$ws.Range("D19").Value = '0.0₃0723'
$ws.Range("E19").Value = '  -1.07%  '
$ws.Range("E20").Value = '  -1.38%  '
$ws.Range("E21").Value = '  -0.09%  '
Set-TextValue $ws "D22" '4.34'
$ws.Range("E22").Value = '  -3.76%  '
Set-TextValue $ws "D23" '9.71'
$ws.Range("E23").Value = '  +3.60%  '
Set-TextValue $ws "D24" '2.03'
$ws.Range("E24").Value = '  -1.24%  '
Set-TextValue $ws "D25" '148.03'
$ws.Range("E25").Value = '  +0.67%  '
$ws.Range("E26").Value = '  -2.83%  '
$ws.Range("E27").Value = '  +1.67%  '
$ws.Range("E30").Value = '  -3.69%  '
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("E32").Value = '  -1.24%  '
Set-TextValue $ws "D33" '3.18'
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("D34").Value = '1.426.06'
$ws.Range("E34").Value = '  -2.68%  '
Set-TextValue $ws "D35" '1.59'
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("E36").Value = '  +0.17%  '
Set-TextValue $ws "D37" '0.568'
$ws.Range("E37").Value = '  -1.56%  '
Set-TextValue $ws "D38" '0.888'
$ws.Range("E38").Value = '  -4.42%  '
$ws.Range("E39").Value = '  -3.28%  '
$ws.Range("E40").Value = '  -1.37%  '
$ws.Range("E41").Value = '  -0.07%  '
Set-TextValue $ws "D42" '0.821'
$ws.Range("E42").Value = '  +3.27%  '
$ws.Range("E43").Value = '  +2.75%  '
$ws.Range("E44").Value = '  -1.82%  '
$ws.Range("E45").Value = '  +0.76%  '
Set-TextValue $ws "D46" '64.90'
$ws.Range("E46").Value = '  -7.24%  '
$ws.Range("D47").Value = '1.788.53'
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("E48").Value = '  -2.71%  '
Set-TextValue $ws "D49" '88.43'
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("E51").Value = '  -3.18%  '
